# Update the "Percent Complete" values in the completeness report.
# Cells B5:B11 hold the Fields-of-Interest percentages, D/E2:E26 hold the
# per-element percentages. The sheet stores these as plain text (shared
# strings), so we set each cell via a text formula and then convert the
# formula result back to a literal value in place (copy / paste-special
# values) -- this avoids Excel's automatic "looks like a number" coercion
# (which would turn the text into a real number) while leaving the cell's
# number format / style completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B5"  = "99.1"
    "B6"  = "97.4"
    "B7"  = "99.1"
    "B8"  = "99.1"
    "B9"  = "97.2"
    "B10" = "96.5"
    "B11" = "91.5"
    "E3"  = "100.0"
    "E4"  = "100.0"
    "E6"  = "100.0"
    "E7"  = "100.0"
    "E8"  = "100.0"
    "E9"  = "100.0"
    "E10" = "100.0"
    "E11" = "100.0"
    "E12" = "100.0"
    "E13" = "93.9"
    "E15" = "31.4"
    "E16" = "24.0"
    "E17" = "19.7"
    "E18" = "16.6"
    "E19" = "19.7"
    "E20" = "22.7"
    "E23" = "100.0"
    "E25" = "100.0"
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $value + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
